# Invoice update: rename "amount_purchased" header to "purchased", and add
# three new purchased-product rows (MacBook Pro M3 MAX, Keychron Mechanical
# Keyboad, Apple Vision Pro) with their quantity / unit price / total price,
# following the same layout/formatting conventions as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "amount_purchased" -> "purchased" ------------------------------
$ws.Range("C1").Value = "purchased"

# --- Extend the bold "product_id" column formatting (used by A2:A4) down
#     into the new rows A5:A7 before writing their values -------------------
$ws.Range("A4").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 5: MacBook Pro M3 MAX -----------------------------------------------
$ws.Range("A5").Value = 9825060
$ws.Range("B5").Value = "MacBook Pro M3 MAX"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4000
$ws.Range("E5").Value = 12000

# --- Row 6: Keychron Mechanical Keyboad -------------------------------------
$ws.Range("A6").Value = 3720602
$ws.Range("B6").Value = "Keychron Mechanical Keyboad"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 400
$ws.Range("E6").Value = 800

# --- Row 7: Apple Vision Pro --------------------------------------------------
$ws.Range("A7").Value = 8690602
$ws.Range("B7").Value = "Apple Vision Pro"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 3500
$ws.Range("E7").Value = 3500

# --- Widen the product_name column to fit the new, longer product names ----
$ws.Columns("B:B").AutoFit()

"done"
